$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview: column G "Latest HO Xliff Generate Date" for rows 4-7
$overview.Range("G4:G7").Value = "2016-08-16 10:31:34"

# zh-cn: column E "Priority" rows 4-7 -> "ht"; column H "Latest Handoff Datetime" rows 4-7
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-16 10:31:28"

# de-de: column E "Priority" rows 4-7 -> "ht"
$dede.Range("E4:E7").Value = "ht"
# de-de: column H "Latest Handoff Datetime" rows 4-7 shares the same underlying
# value as Overview!G4:G7 ("Latest HO Xliff Generate Date"), so it picks up the
# same updated timestamp.
$dede.Range("H4:H7").Value = "2016-08-16 10:31:34"
